# Update 'F' column (想去人数 / interest count) values across sheets
# as described by the commit: 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 347
$ws.Range("F3").Value = 245
$ws.Range("F4").Value = 570
$ws.Range("F5").Value = 1359
$ws.Range("F6").Value = 660
$ws.Range("F7").Value = 355
$ws.Range("F8").Value = 34
$ws.Range("F9").Value = 157
$ws.Range("F10").Value = 419
$ws.Range("F11").Value = 6277
$ws.Range("F12").Value = 117
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 1892
$ws.Range("F15").Value = 4685
$ws.Range("F16").Value = 463
$ws.Range("F19").Value = 5465
$ws.Range("F20").Value = 7190
$ws.Range("F22").Value = 1088
$ws.Range("F23").Value = 758
$ws.Range("F24").Value = 4015
$ws.Range("F25").Value = 560
$ws.Range("F27").Value = 229
$ws.Range("F29").Value = 1063
$ws.Range("F30").Value = 1500
$ws.Range("F31").Value = 557
$ws.Range("F32").Value = 690
$ws.Range("F33").Value = 1690
$ws.Range("F34").Value = 240
$ws.Range("F35").Value = 1894
$ws.Range("F36").Value = 231
$ws.Range("F38").Value = 1246
$ws.Range("F39").Value = 1334
$ws.Range("F40").Value = 692
$ws.Range("F41").Value = 322
$ws.Range("F42").Value = 1168
$ws.Range("F43").Value = 3677
$ws.Range("F44").Value = 158
$ws.Range("F45").Value = 345
$ws.Range("F46").Value = 448
$ws.Range("F47").Value = 24
$ws.Range("F49").Value = 3954

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1271
$ws.Range("F32").Value = 52

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4438

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4438
$ws.Range("F3").Value = 347
$ws.Range("F4").Value = 1271
$ws.Range("F6").Value = 245
$ws.Range("F7").Value = 570
$ws.Range("F9").Value = 1359
$ws.Range("F11").Value = 660
$ws.Range("F12").Value = 355
$ws.Range("F13").Value = 157
$ws.Range("F14").Value = 419
$ws.Range("F15").Value = 117
$ws.Range("F16").Value = 25
$ws.Range("F17").Value = 4685
$ws.Range("F18").Value = 5465
$ws.Range("F19").Value = 5465
$ws.Range("F21").Value = 1088
$ws.Range("F22").Value = 758
$ws.Range("F23").Value = 4015
$ws.Range("F24").Value = 560
$ws.Range("F25").Value = 229
$ws.Range("F28").Value = 1064
$ws.Range("F29").Value = 1500
$ws.Range("F30").Value = 557
$ws.Range("F31").Value = 690
$ws.Range("F32").Value = 1690
$ws.Range("F33").Value = 240
$ws.Range("F34").Value = 1894
$ws.Range("F39").Value = 692
$ws.Range("F41").Value = 322
$ws.Range("F43").Value = 3677
$ws.Range("F45").Value = 158
$ws.Range("F46").Value = 345
$ws.Range("F47").Value = 448
$ws.Range("F50").Value = 3954
$ws.Range("F51").Value = 52
